$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 551
$ws.Range("F4").Value = 577
$ws.Range("F5").Value = 1370
$ws.Range("F6").Value = 679
$ws.Range("F10").Value = 432
$ws.Range("F11").Value = 6363
$ws.Range("F14").Value = 1898
$ws.Range("F15").Value = 4745
$ws.Range("F19").Value = 5543
$ws.Range("F20").Value = 7357
$ws.Range("F23").Value = 765
$ws.Range("F24").Value = 4042
$ws.Range("F25").Value = 566
$ws.Range("F31").Value = 572
$ws.Range("F32").Value = 700
$ws.Range("F33").Value = 1704
$ws.Range("F34").Value = 241
$ws.Range("F35").Value = 1928
$ws.Range("F36").Value = 238
$ws.Range("F37").Value = 46
$ws.Range("F38").Value = 1255
$ws.Range("F40").Value = 699
$ws.Range("F41").Value = 325
$ws.Range("F42").Value = 1625
$ws.Range("F43").Value = 3718
$ws.Range("F44").Value = 161
$ws.Range("F45").Value = 350
$ws.Range("F46").Value = 452
$ws.Range("F48").Value = 101
$ws.Range("F49").Value = 3965

$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 1281
$ws.Range("F17").Value = 19

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 4480

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 4480
$ws.Range("F4").Value = 551
$ws.Range("F5").Value = 1281
$ws.Range("F9").Value = 577
$ws.Range("F11").Value = 1370
$ws.Range("F13").Value = 679
$ws.Range("F17").Value = 432
$ws.Range("F20").Value = 4745
$ws.Range("F21").Value = 5543
$ws.Range("F22").Value = 5543
$ws.Range("F24").Value = 765
$ws.Range("F25").Value = 4042
$ws.Range("F26").Value = 566
$ws.Range("F32").Value = 572
$ws.Range("F33").Value = 700
$ws.Range("F34").Value = 1704
$ws.Range("F35").Value = 241
$ws.Range("F36").Value = 1928
$ws.Range("F41").Value = 699
$ws.Range("F42").Value = 325
$ws.Range("F44").Value = 3718
$ws.Range("F46").Value = 161
$ws.Range("F47").Value = 350
$ws.Range("F48").Value = 101
$ws.Range("F50").Value = 3965
